$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 9 data - "Selection Problem"
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "Selection Algorith"
$ws.Range("C9").Value = "Used to find kth position element in sorted array"
$ws.Range("E9").Value = "O(n2) "
$ws.Range("F9").Value = "O(1)"
$ws.Range("H9").Value = "Selection"

# Update selection to reflect new active cell after edit
$ws.Range("A10").Select()
